$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 1.495631747129835
$arr[1,0] = 1.365099201893429
$arr[2,0] = 1.285384004216667
$arr[3,0] = 1.253008512439408
$arr[4,0] = 1.247639188326445
$arr[5,0] = 1.284946934733341
$arr[6,0] = 1.450534580073054
$arr[7,0] = 1.778683746166791
$arr[8,0] = 2.021898641072028
$arr[9,0] = 2.13301358553241
$arr[10,0] = 2.175158477031346
$arr[11,0] = 2.166078796867737
$arr[12,0] = 2.136479508231446
$arr[13,0] = 2.118357947126299
$arr[14,0] = 2.014646557834908
$arr[15,0] = 1.951144503460625
$arr[16,0] = 1.914664549583392
$arr[17,0] = 1.902320758413111
$arr[18,0] = 1.957899776350871
$arr[19,0] = 2.145171691261112
$arr[20,0] = 2.267961400599916
$arr[21,0] = 2.20239000393957
$arr[22,0] = 1.954845629550903
$arr[23,0] = 1.689540161331706
$ws.Range("B2:B25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.2578220185193629
$arr[1,0] = 0.2424540891787785
$arr[2,0] = 0.2329391894450055
$arr[3,0] = 0.2290420889293898
$arr[4,0] = 0.2283937928183377
$arr[5,0] = 0.2328867113579065
$arr[6,0] = 0.2525396198684859
$arr[7,0] = 0.2904481842083158
$arr[8,0] = 0.3179119590526511
$arr[9,0] = 0.3303213550488522
$arr[10,0] = 0.335008301782608
$arr[11,0] = 0.333999429725452
$arr[12,0] = 0.3307071987203187
$arr[13,0] = 0.328689016369367
$arr[14,0] = 0.3170992730966873
$arr[15,0] = 0.3099677211674532
$arr[16,0] = 0.3058579415454687
$arr[17,0] = 0.3044650884103817
$arr[18,0] = 0.3107277056539317
$arr[19,0] = 0.331674539579268
$arr[20,0] = 0.3452932040565031
$arr[21,0] = 0.3380312360669677
$arr[22,0] = 0.3103841470452267
$arr[23,0] = 0.2802607805170396
$ws.Range("C2:C25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.02687400254731998
$arr[1,0] = 0.02711643910000383
$arr[2,0] = 0.02727541241164388
$arr[3,0] = 0.02734274098959588
$arr[4,0] = 0.02735407466412987
$arr[5,0] = 0.02727631011901721
$arr[6,0] = 0.02695549600567304
$arr[7,0] = 0.02640658870086554
$arr[8,0] = 0.02605215394381233
$arr[9,0] = 0.02590151975193677
$arr[10,0] = 0.0258460039156958
$arr[11,0] = 0.02585789233679492
$arr[12,0] = 0.02589692183437364
$arr[13,0] = 0.02592102732922186
$arr[14,0] = 0.02606221168206346
$arr[15,0] = 0.02615153972320527
$arr[16,0] = 0.0262039164256489
$arr[17,0] = 0.02622182156591002
$arr[18,0] = 0.02614192734509757
$arr[19,0] = 0.02588541649075538
$arr[20,0] = 0.02572666810060653
$arr[21,0] = 0.02581058035689665
$arr[22,0] = 0.02614626992107816
$arr[23,0] = 0.02654650298055294
$ws.Range("D2:D25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 1.047172809179429
$arr[1,0] = 1.03165498405393
$arr[2,0] = 1.023008281313622
$arr[3,0] = 1.019705117991904
$arr[4,0] = 1.019169909881882
$arr[5,0] = 1.02296284271219
$arr[6,0] = 1.041638634564208
$arr[7,0] = 1.085310790662845
$arr[8,0] = 1.121776611643
$arr[9,0] = 1.13933502951015
$arr[10,0] = 1.146124794757426
$arr[11,0] = 1.144656215663488
$arr[12,0] = 1.139890798820474
$arr[13,0] = 1.136990214877031
$arr[14,0] = 1.120648756149563
$arr[15,0] = 1.110873145935571
$arr[16,0] = 1.105341690322689
$arr[17,0] = 1.103484461579896
$arr[18,0] = 1.111904326367238
$arr[19,0] = 1.141286686340734
$arr[20,0] = 1.161310898563826
$arr[21,0] = 1.1505480349267
$arr[22,0] = 1.111437853585571
$arr[23,0] = 1.072733144581917
$ws.Range("F2:F25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.9128616925189732
$arr[1,0] = 0.8970943396398212
$arr[2,0] = 0.8882699481462026
$arr[3,0] = 0.8848879585360834
$arr[4,0] = 0.8843392661573262
$arr[5,0] = 0.8882234727912106
$arr[6,0] = 0.9072464043021142
$arr[7,0] = 0.9514166244728131
$arr[8,0] = 0.9881526197793846
$arr[9,0] = 1.005816481374524
$arr[10,0] = 1.012643939424521
$arr[11,0] = 1.011167337611113
$arr[12,0] = 1.006375394944996
$arr[13,0] = 1.003458280058055
$arr[14,0] = 0.9870175453972365
$arr[15,0] = 0.9771767558577551
$arr[16,0] = 0.9716061510154361
$arr[17,0] = 0.9697353735867864
$arr[18,0] = 0.9782150445325328
$arr[19,0] = 1.007779132253376
$arr[20,0] = 1.027909205632852
$arr[21,0] = 1.017090932399725
$arr[22,0] = 0.977745363360043
$arr[23,0] = 0.9387221814268401
$ws.Range("G2:G25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.9087555831881531
$arr[1,0] = 0.9079242766196955
$arr[2,0] = 0.9080481364246964
$arr[3,0] = 0.9082575760753429
$arr[4,0] = 0.9083019390140947
$arr[5,0] = 0.9080503180661879
$arr[6,0] = 0.9083369803885262
$arr[7,0] = 0.9139571724169002
$arr[8,0] = 0.9212068031884257
$arr[9,0] = 0.9251905465093273
$arr[10,0] = 0.9267983378686893
$arr[11,0] = 0.9264476480169037
$arr[12,0] = 0.9253208281671732
$arr[13,0] = 0.9246435598485334
$arr[14,0] = 0.9209603093079295
$arr[15,0] = 0.9188768499564333
$arr[16,0] = 0.9177430180635895
$arr[17,0] = 0.9173701846950451
$arr[18,0] = 0.9190919559237898
$arr[19,0] = 0.9256491044068298
$arr[20,0] = 0.9305132782644421
$arr[21,0] = 0.9278640238734113
$arr[22,0] = 0.9189945072772332
$arr[23,0] = 0.9118914524562172
$ws.Range("H2:H25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.3022342068966566
$arr[1,0] = 0.2910530132305098
$arr[2,0] = 0.2843731862990069
$arr[3,0] = 0.2816975864579234
$arr[4,0] = 0.2812561082109255
$arr[5,0] = 0.2843369142118632
$arr[6,0] = 0.2983403389382602
$arr[7,0] = 0.3272814963537769
$arr[8,0] = 0.3494620940073503
$arr[9,0] = 0.3597552948515528
$arr[10,0] = 0.3636825006458935
$arr[11,0] = 0.3628353963484017
$arr[12,0] = 0.3600777982728545
$arr[13,0] = 0.35839252329896
$arr[14,0] = 0.3487935123307295
$arr[15,0] = 0.3429569873279092
$arr[16,0] = 0.33961908144191
$arr[17,0] = 0.3384922003073001
$arr[18,0] = 0.3435763163565611
$arr[19,0] = 0.3608869724494639
$arr[20,0] = 0.3723719053383547
$arr[21,0] = 0.3662264342279826
$arr[22,0] = 0.3432962629290728
$arr[23,0] = 0.3192921134006497
$ws.Range("L2:L25").Value = $arr

